# #5: property aircraft done
#
# - The "建物" (building) sheet had its property_category column (I2:I17)
#   mistakenly copied from the "land" sheet as "land" instead of "building".
#   Fix every data row to read "building".
# - The "航空器" (aircraft) sheet is removed entirely from the workbook.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Fix the mislabeled property_category values on the "建物" (building) sheet.
$wsBuilding = $wb.Worksheets.Item("建物")
for ($row = 2; $row -le 17; $row++) {
    $wsBuilding.Cells.Item($row, 9).Value = "building"
}

# Remove the "航空器" (aircraft) sheet completely.
$wsAircraft = $wb.Worksheets.Item("航空器")
$wsAircraft.Delete()
